$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two stale data columns (1989, 2002) - the "2014" column (D) shifts left into B
$ws.Range("B1:C1").EntireColumn.Delete()

# Move the title text down onto the subtitle row (which carries the style the title should
# use), drop the stray "census results" subtitle text, then remove the now-empty old title
# row so everything shifts up into place.
$ws.Range("A2").Value = $ws.Range("A1").Value()
$ws.Rows.Item(1).Delete()

# Match the sheet tab name to the municipality
$ws.Name = "ვანი"

$ws.Range("A2").Select()
